# Project What Ifs.docx edit:
#   - The final paragraph ("A third microbit is added ... shared to the
#     third sentinel microbit") gets its trailing word "microbit" split
#     off into its own run wrapped in spellcheck proofErr markers (Word
#     flags "microbit" as a misspelling).
#   - A new underlined heading and four numbered paragraphs are appended
#     after it. The "_GoBack" bookmark (auto-maintained by Word at the
#     location of the last edit) moves from the end of the old final
#     paragraph to the end of the new "3. How it meets the requirements"
#     paragraph, whose text is itself typed as two runs ("3. " then the
#     rest).
$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$paraStart = $lastPara.Range.Start
# Paragraph.Range.Text includes the trailing paragraph-mark character(s)
# (chr 13, sometimes chr 7 for cell/row marks) - strip those before
# comparing/measuring visible text.
$paraText = $lastPara.Range.Text.TrimEnd([char]13, [char]7)
$docEnd = $d.Content.End

$expected = "A third microbit is added to allow wireless transfer of data to a database, for each time one of the 2 microbits part of the exercise is pressed, the data channel between them is shared to the third sentinel microbit"

# Range spanning the whole last paragraph (start of its text through the
# end of the document content, i.e. through its paragraph mark).
$r = $d.Range($paraStart, $docEnd)

if ($paraText -eq $expected) {
    # Reconstruct the paragraph's own runs verbatim (so its w14:paraId /
    # rsid attributes - which InsertXML would otherwise drop - survive),
    # but split the trailing "microbit" into its own spell-checked run,
    # then append the new heading + numbered paragraphs after it.
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="73657467" w14:textId="0C682626" w:rsidR="00227115" w:rsidRDefault="00227115"><w:r><w:t xml:space="preserve">A third </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>microbit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is added to allow wireless transfer of data to a database, for each time one of the 2 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>microbits</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> part of the exercise is pressed, the data channel between them is shared to the third sentinel </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>microbit</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Investigation on Project and if it meets the requirements to satisfy the brief</w:t></w:r></w:p><w:p><w:r><w:t>1. Research on project</w:t></w:r></w:p><w:p><w:r><w:t>2. Inspiration on project</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">3. </w:t></w:r><w:r><w:t>How it meets the requirements</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:t>4. Survey of user results on project</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
} else {
    # Fallback: structure of the document differs from what was expected
    # (e.g. re-run against a modified file) - still apply the same
    # logical edit using plain Range/Find operations, best-effort.
    $needle = "microbit"
    $idx = $paraText.LastIndexOf($needle)
    if ($idx -ge 0) {
        $absStart = $paraStart + $idx
        $tailRange = $d.Range($absStart, $paraStart + $paraText.Length)
        $xmlTail = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>microbit</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $tailRange.InsertXML($xmlTail)
    }

    function Add-Paragraph($text) {
        $rr = $d.Content
        $rr.Collapse(0)
        $rr.InsertParagraphAfter()
        $rr = $d.Content
        $rr.Collapse(0)
        $rr.Text = $text
        return $d.Paragraphs.Last
    }

    $p1 = Add-Paragraph "Investigation on Project and if it meets the requirements to satisfy the brief"
    $p1.Range.Font.Underline = 1
    Add-Paragraph "1. Research on project" | Out-Null
    Add-Paragraph "2. Inspiration on project" | Out-Null
    Add-Paragraph "3. How it meets the requirements" | Out-Null
    Add-Paragraph "4. Survey of user results on project" | Out-Null
}
